$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Correct the "Solar" column (E) values: the values were stored in
#        Watts (x1000) but should be in Kilowatts. Divide the affected
#        cells down to kW. ---
$ws.Range("E13").Value = 21
$ws.Range("E19").Value = 11.4
$ws.Range("E21").Value = 7.6
$ws.Range("E22").Value = 19.3
$ws.Range("E23").Value = 31.9
$ws.Range("E24").Value = 9.279999999999999
$ws.Range("E25").Value = 6
$ws.Range("E26").Value = 13.96

# --- 2. The shared custom number format used by the data cells (B2:G26)
#        needs one more decimal place so fractional kW values show. ---
$ws.Range("B2:G26").NumberFormat = "#,##0.0"

# --- 3. Update the chart: axis title text and the value-axis number
#        format, reverting back to plain kilowatts. ---
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$valueAxis = $chart.Axes(2)
$valueAxis.AxisTitle.Text = "Kilowatts (kW)"
$valueAxis.NumberFormat = "#,##0"
